$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 387, shifting the existing rows 387:496 down to 388:497
$ws.Rows(387).Insert()

# Populate the newly inserted row 387 with its data
$ws.Range("A387").Value = 5
$ws.Range("B387").Value = "Macroferia Regional de Talca"
$ws.Range("C387").Value = "Maule"
$ws.Range("D387").Value = 44841
$ws.Range("E387").Value = 7
$ws.Range("F387").Value = 100112043
$ws.Range("G387").Value = "Pepino ensalada"
$ws.Range("H387").Value = "Sin especificar"
$ws.Range("I387").Value = "Primera"
$ws.Range("J387").Value = 300
$ws.Range("K387").Value = 25000
$ws.Range("L387").Value = 25000
$ws.Range("M387").Value = 25000
$ws.Range("N387").Value = "$/caja 60 unidades"
$ws.Range("O387").Value = "Región de Arica y Parinacota"
$ws.Range("P387").Value = 417
$ws.Range("Q387").Value = 60
$ws.Range("R387").Value = "Hortaliza"

# Match the date-number-format style used by the other cells in column D
$ws.Range("D387").NumberFormat = $ws.Range("D388").NumberFormat
